# Applies updated crypto price/volume data to Sheet1 (cryptos.xlsx)
# Matches commit "Updated cryptos list ... with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '38.689.00'
$c.Style = "Normal"
$ws.Range("E2").Value = '  +2.58%  '

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '2.085.30'
$c.Style = "Normal"
$ws.Range("E3").Value = '  +1.60%  '

$ws.Range("E4").Value = '  +0.00%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '228.34'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +0.32%  '

$ws.Range("E6").Value = '  +0.81%  '

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '59.96'
$c.Style = "Normal"
$ws.Range("E7").Value = '  +1.24%  '

$ws.Range("E8").Value = '  -0.03%  '

$ws.Range("E9").Value = '  +2.09%  '

$ws.Range("E10").Value = '  +0.50%  '

$ws.Range("E11").Value = '  -0.34%  '

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '2.394.20'
$c.Style = "Normal"
$ws.Range("E12").Value = '  +1.57%  '

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '15.03'
$c.Style = "Normal"
$ws.Range("E13").Value = '  +4.13%  '

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '21.90'
$c.Style = "Normal"
$ws.Range("E14").Value = '  +2.30%  '

$ws.Range("E15").Value = '  +4.76%  '

$ws.Range("E16").Value = '  +0.09%  '

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '2.084.73'
$c.Style = "Normal"
$ws.Range("E17").Value = '  +1.26%  '

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '38.634.79'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +2.44%  '

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '71.49'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +2.95%  '

$ws.Range("E20").Value = '  +1.03%  '

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '0.0₃0838'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +0.77%  '

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '226.92'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +1.93%  '

$ws.Range("E23").Value = '  -0.45%  '

$ws.Range("E24").Value = '  -1.88%  '

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '2.34'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +2.58%  '

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '171.01'
$c.Style = "Normal"
$ws.Range("E26").Value = '  +1.11%  '

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '9.54'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +2.55%  '

$ws.Range("E28").Value = '  +7.06%  '

$ws.Range("E29").Value = '  +13.31%  '

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '19.15'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +1.90%  '

$ws.Range("E31").Value = '  +0.99%  '

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '2.39'
$c.Style = "Normal"
$ws.Range("E32").Value = '  +4.61%  '

$ws.Range("E33").Value = '  +2.52%  '

$ws.Range("E34").Value = '  +2.90%  '

$ws.Range("E35").Value = '  +0.85%  '

$ws.Range("E36").Value = '  -0.04%  '

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '2.39'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +1.28%  '

$ws.Range("E38").Value = '  +3.30%  '

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -0.13%  '

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '17.90'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -2.62%  '

$ws.Range("E41").Value = '  +5.72%  '

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '1.544.93'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +1.16%  '

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '100.21'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +1.83%  '

$ws.Range("E44").Value = '  +3.62%  '

$ws.Range("E45").Value = '  -1.00%  '

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '7.72'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +8.98%  '

$ws.Range("E47").Value = '  +0.70%  '

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '4.11'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -0.60%  '

$ws.Range("E49").Value = '  +2.62%  '

$ws.Range("E50").Value = '  +0.25%  '

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '2.283.02'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +1.69%  '
